$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: additional time logged against "today" (crouching and aim offset work)
$ws.Range("C34").Formula = "=(1/60)*(7+20)"
$ws.Range("D34").Formula = "=(1/60)*(22+10+20+20+20+20+5+20+5+20)"

# Remaining-days estimates now rounded up to whole days
$ws.Range("I7").Formula = "=ROUNDUP(I3/I5, 0)"
$ws.Range("I9").Formula = "=ROUNDUP(I3/F34, 0)"

# I8 used to read "BASED ON LAST DAY"; it now reflects "today's" rate,
# and a new row (I10/I11) is added to keep the old "last day" estimate
$ws.Range("I8").Value = "BASED ON TODAY"

$ws.Range("I10").Value = "BASED ON LAST DAY"
$ws.Range("I10").Font.Bold = $true
$ws.Range("I10").NumberFormat = "0.00"

$ws.Range("I11").Formula = "=ROUNDUP(I3/F33, 0)"
